# Auto-generated PowerShell Excel COM-interop edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update reordered rows (content moved between rows sharing the same date-slot) ---
# Row 2
$ws.Cells.Item(2, 6).Value = 'Zambrow'
$ws.Cells.Item(2, 7).Value = 1
$ws.Cells.Item(2, 8).Value = 'Legia II'
$ws.Cells.Item(2, 9).Value = 4
$ws.Cells.Item(2, 10).Value = 3.32
$ws.Cells.Item(2, 11).Value = '04/08/2023 00:12'
$ws.Cells.Item(2, 12).Value = 3.47
$ws.Cells.Item(2, 13).Value = '05/08/2023 11:58'
$ws.Cells.Item(2, 14).Value = 3.52
$ws.Cells.Item(2, 15).Value = '04/08/2023 00:12'
$ws.Cells.Item(2, 16).Value = 3.7
$ws.Cells.Item(2, 17).Value = '05/08/2023 11:58'
$ws.Cells.Item(2, 18).Value = 1.78
$ws.Cells.Item(2, 19).Value = '04/08/2023 00:12'
$ws.Cells.Item(2, 20).Value = 1.83
$ws.Cells.Item(2, 21).Value = '05/08/2023 11:58'
$ws.Cells.Item(2, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/olimpia-zambrow-legia/jkkoUWkl/'

# Row 3
$ws.Cells.Item(3, 6).Value = 'Jagiellonia II'
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(3, 8).Value = 'Bron Radom'
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1.63
$ws.Cells.Item(3, 11).Value = '04/08/2023 00:12'
$ws.Cells.Item(3, 12).Value = 1.17
$ws.Cells.Item(3, 13).Value = '05/08/2023 11:58'
$ws.Cells.Item(3, 14).Value = 3.81
$ws.Cells.Item(3, 15).Value = '04/08/2023 00:12'
$ws.Cells.Item(3, 16).Value = 7.01
$ws.Cells.Item(3, 17).Value = '05/08/2023 11:58'
$ws.Cells.Item(3, 18).Value = 3.63
$ws.Cells.Item(3, 19).Value = '04/08/2023 00:12'
$ws.Cells.Item(3, 20).Value = 8.640000000000001
$ws.Cells.Item(3, 21).Value = '05/08/2023 11:58'
$ws.Cells.Item(3, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/jagiellonia-bron-radom/K2SDnu6F/'

# Row 58
$ws.Cells.Item(58, 6).Value = 'T. Mazowiecki'
$ws.Cells.Item(58, 7).Value = 4
$ws.Cells.Item(58, 8).Value = 'Concordia Elblag'
$ws.Cells.Item(58, 9).Value = 1
$ws.Cells.Item(58, 10).Value = 1.46
$ws.Cells.Item(58, 11).Value = '15/09/2023 03:13'
$ws.Cells.Item(58, 12).Value = 1.49
$ws.Cells.Item(58, 13).Value = '16/09/2023 15:36'
$ws.Cells.Item(58, 14).Value = 3.97
$ws.Cells.Item(58, 15).Value = '15/09/2023 03:13'
$ws.Cells.Item(58, 16).Value = 4.17
$ws.Cells.Item(58, 17).Value = '16/09/2023 15:36'
$ws.Cells.Item(58, 18).Value = 4.65
$ws.Cells.Item(58, 19).Value = '15/09/2023 03:13'
$ws.Cells.Item(58, 20).Value = 5.16
$ws.Cells.Item(58, 21).Value = '16/09/2023 15:36'
$ws.Cells.Item(58, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/tomaszow-mazowiecki-concordia-elblag/EJtQZeN2/'

# Row 59
$ws.Cells.Item(59, 6).Value = 'Swit Mazowiecki'
$ws.Cells.Item(59, 7).Value = 1
$ws.Cells.Item(59, 8).Value = 'Pelikan'
$ws.Cells.Item(59, 9).Value = 1
$ws.Cells.Item(59, 10).Value = 1.69
$ws.Cells.Item(59, 11).Value = '15/09/2023 03:13'
$ws.Cells.Item(59, 12).Value = 1.72
$ws.Cells.Item(59, 13).Value = '16/09/2023 15:45'
$ws.Cells.Item(59, 14).Value = 3.5
$ws.Cells.Item(59, 15).Value = '15/09/2023 03:13'
$ws.Cells.Item(59, 16).Value = 3.9
$ws.Cells.Item(59, 17).Value = '16/09/2023 15:45'
$ws.Cells.Item(59, 18).Value = 3.65
$ws.Cells.Item(59, 19).Value = '15/09/2023 03:13'
$ws.Cells.Item(59, 20).Value = 3.74
$ws.Cells.Item(59, 21).Value = '16/09/2023 15:45'
$ws.Cells.Item(59, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/swit-mazowiecki-pelikan/zZUYXZiF/'

# Row 60
$ws.Cells.Item(60, 6).Value = 'Pilica Bialobrzegi'
$ws.Cells.Item(60, 7).Value = 1
$ws.Cells.Item(60, 8).Value = 'Legia II'
$ws.Cells.Item(60, 9).Value = 1
$ws.Cells.Item(60, 10).Value = 3.93
$ws.Cells.Item(60, 11).Value = '15/09/2023 03:13'
$ws.Cells.Item(60, 12).Value = 6.37
$ws.Cells.Item(60, 13).Value = '16/09/2023 15:56'
$ws.Cells.Item(60, 14).Value = 3.78
$ws.Cells.Item(60, 15).Value = '15/09/2023 03:13'
$ws.Cells.Item(60, 16).Value = 4.7
$ws.Cells.Item(60, 17).Value = '16/09/2023 15:56'
$ws.Cells.Item(60, 18).Value = 1.61
$ws.Cells.Item(60, 19).Value = '15/09/2023 03:13'
$ws.Cells.Item(60, 20).Value = 1.36
$ws.Cells.Item(60, 21).Value = '16/09/2023 15:56'
$ws.Cells.Item(60, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/pilica-bialobrzegi-legia/lxZxXg7L/'

# Row 61
$ws.Cells.Item(61, 6).Value = 'Mlawa'
$ws.Cells.Item(61, 7).Value = 3
$ws.Cells.Item(61, 8).Value = 'Bron Radom'
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 1.76
$ws.Cells.Item(61, 11).Value = '15/09/2023 03:13'
$ws.Cells.Item(61, 12).Value = 1.95
$ws.Cells.Item(61, 13).Value = '16/09/2023 15:46'
$ws.Cells.Item(61, 14).Value = 3.76
$ws.Cells.Item(61, 15).Value = '15/09/2023 03:13'
$ws.Cells.Item(61, 16).Value = 3.76
$ws.Cells.Item(61, 17).Value = '16/09/2023 15:46'
$ws.Cells.Item(61, 18).Value = 3.16
$ws.Cells.Item(61, 19).Value = '15/09/2023 03:13'
$ws.Cells.Item(61, 20).Value = 3.05
$ws.Cells.Item(61, 21).Value = '16/09/2023 15:46'
$ws.Cells.Item(61, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/mks-mlawa-bron-radom/hAuUYFx9/'

# Row 62
$ws.Cells.Item(62, 6).Value = 'LKS Lomza'
$ws.Cells.Item(62, 7).Value = 2
$ws.Cells.Item(62, 8).Value = 'Jagiellonia II'
$ws.Cells.Item(62, 9).Value = 1
$ws.Cells.Item(62, 10).Value = 2.62
$ws.Cells.Item(62, 11).Value = '16/09/2023 12:43'
$ws.Cells.Item(62, 12).Value = 2.73
$ws.Cells.Item(62, 13).Value = '16/09/2023 15:58'
$ws.Cells.Item(62, 14).Value = 3.3
$ws.Cells.Item(62, 15).Value = '16/09/2023 12:43'
$ws.Cells.Item(62, 16).Value = 3.62
$ws.Cells.Item(62, 17).Value = '16/09/2023 15:58'
$ws.Cells.Item(62, 18).Value = 2.29
$ws.Cells.Item(62, 19).Value = '16/09/2023 12:43'
$ws.Cells.Item(62, 20).Value = 2.16
$ws.Cells.Item(62, 21).Value = '16/09/2023 15:58'
$ws.Cells.Item(62, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/lks-lomza-jagiellonia/2qFDQipk/'

# Row 63
$ws.Cells.Item(63, 6).Value = 'Warta Sieradz'
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 'Wikielec'
$ws.Cells.Item(63, 9).Value = 1
$ws.Cells.Item(63, 10).Value = 2.4
$ws.Cells.Item(63, 11).Value = '15/09/2023 03:13'
$ws.Cells.Item(63, 12).Value = 2.61
$ws.Cells.Item(63, 13).Value = '16/09/2023 15:59'
$ws.Cells.Item(63, 14).Value = 3.18
$ws.Cells.Item(63, 15).Value = '15/09/2023 03:13'
$ws.Cells.Item(63, 16).Value = 3.47
$ws.Cells.Item(63, 17).Value = '16/09/2023 15:59'
$ws.Cells.Item(63, 18).Value = 2.4
$ws.Cells.Item(63, 19).Value = '15/09/2023 03:13'
$ws.Cells.Item(63, 20).Value = 2.31
$ws.Cells.Item(63, 21).Value = '16/09/2023 15:59'
$ws.Cells.Item(63, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/warta-sieradz-gks-wikielec/zcDLOVE1/'

# Row 68
$ws.Cells.Item(68, 6).Value = 'Wikielec'
$ws.Cells.Item(68, 7).Value = 1
$ws.Cells.Item(68, 8).Value = 'T. Mazowiecki'
$ws.Cells.Item(68, 9).Value = 1
$ws.Cells.Item(68, 10).Value = 2.22
$ws.Cells.Item(68, 11).Value = '22/09/2023 02:13'
$ws.Cells.Item(68, 12).Value = 2.5
$ws.Cells.Item(68, 13).Value = '23/09/2023 14:49'
$ws.Cells.Item(68, 14).Value = 3.21
$ws.Cells.Item(68, 15).Value = '22/09/2023 02:13'
$ws.Cells.Item(68, 16).Value = 3.48
$ws.Cells.Item(68, 17).Value = '23/09/2023 14:49'
$ws.Cells.Item(68, 18).Value = 2.67
$ws.Cells.Item(68, 19).Value = '22/09/2023 02:13'
$ws.Cells.Item(68, 20).Value = 2.39
$ws.Cells.Item(68, 21).Value = '23/09/2023 14:49'
$ws.Cells.Item(68, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/gks-wikielec-tomaszow-mazowiecki/nXGPNkU7/'

# Row 69
$ws.Cells.Item(69, 6).Value = 'Jagiellonia II'
$ws.Cells.Item(69, 7).Value = 2
$ws.Cells.Item(69, 8).Value = 'Sulejowek'
$ws.Cells.Item(69, 9).Value = 2
$ws.Cells.Item(69, 10).Value = 2.28
$ws.Cells.Item(69, 11).Value = '22/09/2023 02:13'
$ws.Cells.Item(69, 12).Value = 2.78
$ws.Cells.Item(69, 13).Value = '23/09/2023 14:46'
$ws.Cells.Item(69, 14).Value = 3.28
$ws.Cells.Item(69, 15).Value = '22/09/2023 02:13'
$ws.Cells.Item(69, 16).Value = 3.35
$ws.Cells.Item(69, 17).Value = '23/09/2023 14:50'
$ws.Cells.Item(69, 18).Value = 2.48
$ws.Cells.Item(69, 19).Value = '22/09/2023 02:13'
$ws.Cells.Item(69, 20).Value = 2.15
$ws.Cells.Item(69, 21).Value = '23/09/2023 14:46'
$ws.Cells.Item(69, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/jagiellonia-victoria-sulejowek/SSgCogME/'

# Row 70
$ws.Cells.Item(70, 6).Value = 'Zambrow'
$ws.Cells.Item(70, 7).Value = 3
$ws.Cells.Item(70, 8).Value = 'Grodzisk M.'
$ws.Cells.Item(70, 9).Value = 3
$ws.Cells.Item(70, 10).Value = 4.96
$ws.Cells.Item(70, 11).Value = '22/09/2023 03:13'
$ws.Cells.Item(70, 12).Value = 3.47
$ws.Cells.Item(70, 13).Value = '23/09/2023 15:42'
$ws.Cells.Item(70, 14).Value = 4.12
$ws.Cells.Item(70, 15).Value = '22/09/2023 03:13'
$ws.Cells.Item(70, 16).Value = 3.69
$ws.Cells.Item(70, 17).Value = '23/09/2023 15:42'
$ws.Cells.Item(70, 18).Value = 1.43
$ws.Cells.Item(70, 19).Value = '22/09/2023 03:13'
$ws.Cells.Item(70, 20).Value = 1.83
$ws.Cells.Item(70, 21).Value = '23/09/2023 15:42'
$ws.Cells.Item(70, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/olimpia-zambrow-gks-pogon-grodzisk-mazowiecki/zgcKqXjR/'

# Row 71
$ws.Cells.Item(71, 6).Value = 'Legionowo'
$ws.Cells.Item(71, 7).Value = 1
$ws.Cells.Item(71, 8).Value = 'Warta Sieradz'
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 1.6
$ws.Cells.Item(71, 11).Value = '22/09/2023 03:13'
$ws.Cells.Item(71, 12).Value = 1.77
$ws.Cells.Item(71, 13).Value = '23/09/2023 15:55'
$ws.Cells.Item(71, 14).Value = 3.75
$ws.Cells.Item(71, 15).Value = '22/09/2023 03:13'
$ws.Cells.Item(71, 16).Value = 3.78
$ws.Cells.Item(71, 17).Value = '23/09/2023 15:55'
$ws.Cells.Item(71, 18).Value = 3.85
$ws.Cells.Item(71, 19).Value = '22/09/2023 03:13'
$ws.Cells.Item(71, 20).Value = 3.63
$ws.Cells.Item(71, 21).Value = '23/09/2023 15:55'
$ws.Cells.Item(71, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/legionowo-warta-sieradz/lxe8nZ68/'

# Row 76
$ws.Cells.Item(76, 6).Value = 'Grodzisk M.'
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = 'Legia II'
$ws.Cells.Item(76, 9).Value = 1
$ws.Cells.Item(76, 10).Value = 2.03
$ws.Cells.Item(76, 11).Value = '29/09/2023 02:13'
$ws.Cells.Item(76, 12).Value = 2.03
$ws.Cells.Item(76, 13).Value = '30/09/2023 14:17'
$ws.Cells.Item(76, 14).Value = 3.34
$ws.Cells.Item(76, 15).Value = '29/09/2023 02:13'
$ws.Cells.Item(76, 16).Value = 3.51
$ws.Cells.Item(76, 17).Value = '30/09/2023 14:17'
$ws.Cells.Item(76, 18).Value = 2.82
$ws.Cells.Item(76, 19).Value = '29/09/2023 02:13'
$ws.Cells.Item(76, 20).Value = 3.04
$ws.Cells.Item(76, 21).Value = '30/09/2023 14:17'
$ws.Cells.Item(76, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/gks-pogon-grodzisk-mazowiecki-legia/rZljcjLQ/'

# Row 77
$ws.Cells.Item(77, 6).Value = 'Pilica Bialobrzegi'
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = 'Bron Radom'
$ws.Cells.Item(77, 9).Value = 4
$ws.Cells.Item(77, 10).Value = 1.97
$ws.Cells.Item(77, 11).Value = '29/09/2023 02:12'
$ws.Cells.Item(77, 12).Value = 2.25
$ws.Cells.Item(77, 13).Value = '30/09/2023 14:59'
$ws.Cells.Item(77, 14).Value = 3.36
$ws.Cells.Item(77, 15).Value = '29/09/2023 02:12'
$ws.Cells.Item(77, 16).Value = 3.71
$ws.Cells.Item(77, 17).Value = '30/09/2023 14:59'
$ws.Cells.Item(77, 18).Value = 2.88
$ws.Cells.Item(77, 19).Value = '29/09/2023 02:12'
$ws.Cells.Item(77, 20).Value = 2.56
$ws.Cells.Item(77, 21).Value = '30/09/2023 14:59'
$ws.Cells.Item(77, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/pilica-bialobrzegi-bron-radom/n7esaCjE/'

# Row 78
$ws.Cells.Item(78, 6).Value = 'Warta Sieradz'
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 'Jagiellonia II'
$ws.Cells.Item(78, 9).Value = 1
$ws.Cells.Item(78, 10).Value = 2.36
$ws.Cells.Item(78, 11).Value = '29/09/2023 02:13'
$ws.Cells.Item(78, 12).Value = 2.22
$ws.Cells.Item(78, 13).Value = '30/09/2023 14:50'
$ws.Cells.Item(78, 14).Value = 3.27
$ws.Cells.Item(78, 15).Value = '29/09/2023 02:13'
$ws.Cells.Item(78, 16).Value = 3.4
$ws.Cells.Item(78, 17).Value = '30/09/2023 13:50'
$ws.Cells.Item(78, 18).Value = 2.39
$ws.Cells.Item(78, 19).Value = '29/09/2023 02:13'
$ws.Cells.Item(78, 20).Value = 2.77
$ws.Cells.Item(78, 21).Value = '30/09/2023 14:50'
$ws.Cells.Item(78, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/warta-sieradz-jagiellonia/h8vQjlbf/'

# Row 79
$ws.Cells.Item(79, 6).Value = 'Swit Mazowiecki'
$ws.Cells.Item(79, 7).Value = 6
$ws.Cells.Item(79, 8).Value = 'Concordia Elblag'
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 1.35
$ws.Cells.Item(79, 11).Value = '29/09/2023 03:13'
$ws.Cells.Item(79, 12).Value = 1.36
$ws.Cells.Item(79, 13).Value = '30/09/2023 15:43'
$ws.Cells.Item(79, 14).Value = 4.27
$ws.Cells.Item(79, 15).Value = '29/09/2023 03:13'
$ws.Cells.Item(79, 16).Value = 4.92
$ws.Cells.Item(79, 17).Value = '30/09/2023 15:43'
$ws.Cells.Item(79, 18).Value = 5.68
$ws.Cells.Item(79, 19).Value = '29/09/2023 03:13'
$ws.Cells.Item(79, 20).Value = 6.02
$ws.Cells.Item(79, 21).Value = '30/09/2023 15:43'
$ws.Cells.Item(79, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/swit-mazowiecki-concordia-elblag/Y37azTDD/'

# Row 80
$ws.Cells.Item(80, 6).Value = 'T. Mazowiecki'
$ws.Cells.Item(80, 7).Value = 2
$ws.Cells.Item(80, 8).Value = 'Mlawa'
$ws.Cells.Item(80, 9).Value = 1
$ws.Cells.Item(80, 10).Value = 1.55
$ws.Cells.Item(80, 11).Value = '29/09/2023 03:13'
$ws.Cells.Item(80, 12).Value = 1.34
$ws.Cells.Item(80, 13).Value = '30/09/2023 15:58'
$ws.Cells.Item(80, 14).Value = 3.99
$ws.Cells.Item(80, 15).Value = '29/09/2023 03:13'
$ws.Cells.Item(80, 16).Value = 5.01
$ws.Cells.Item(80, 17).Value = '30/09/2023 15:58'
$ws.Cells.Item(80, 18).Value = 3.93
$ws.Cells.Item(80, 19).Value = '29/09/2023 03:13'
$ws.Cells.Item(80, 20).Value = 6.34
$ws.Cells.Item(80, 21).Value = '30/09/2023 15:58'
$ws.Cells.Item(80, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/tomaszow-mazowiecki-mks-mlawa/fD8ey9b7/'

# Row 87
$ws.Cells.Item(87, 6).Value = 'Concordia Elblag'
$ws.Cells.Item(87, 7).Value = 2
$ws.Cells.Item(87, 8).Value = 'Pilica Bialobrzegi'
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 2.03
$ws.Cells.Item(87, 11).Value = '06/10/2023 02:12'
$ws.Cells.Item(87, 12).Value = 1.83
$ws.Cells.Item(87, 13).Value = '06/10/2023 03:13'
$ws.Cells.Item(87, 14).Value = 3.35
$ws.Cells.Item(87, 15).Value = '06/10/2023 02:12'
$ws.Cells.Item(87, 16).Value = 3.54
$ws.Cells.Item(87, 17).Value = '07/10/2023 13:05'
$ws.Cells.Item(87, 18).Value = 2.81
$ws.Cells.Item(87, 19).Value = '06/10/2023 02:12'
$ws.Cells.Item(87, 20).Value = 3.33
$ws.Cells.Item(87, 21).Value = '06/10/2023 03:13'
$ws.Cells.Item(87, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/concordia-elblag-pilica-bialobrzegi/QJH3yUF4/'

# Row 89
$ws.Cells.Item(89, 6).Value = 'Zambrow'
$ws.Cells.Item(89, 7).Value = 3
$ws.Cells.Item(89, 8).Value = 'Sulejowek'
$ws.Cells.Item(89, 9).Value = 1
$ws.Cells.Item(89, 10).Value = 2.88
$ws.Cells.Item(89, 11).Value = '06/10/2023 02:12'
$ws.Cells.Item(89, 12).Value = 2.56
$ws.Cells.Item(89, 13).Value = '07/10/2023 12:40'
$ws.Cells.Item(89, 14).Value = 3.27
$ws.Cells.Item(89, 15).Value = '06/10/2023 02:12'
$ws.Cells.Item(89, 16).Value = 3.32
$ws.Cells.Item(89, 17).Value = '07/10/2023 13:02'
$ws.Cells.Item(89, 18).Value = 2.02
$ws.Cells.Item(89, 19).Value = '06/10/2023 02:12'
$ws.Cells.Item(89, 20).Value = 2.4
$ws.Cells.Item(89, 21).Value = '07/10/2023 12:40'
$ws.Cells.Item(89, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/olimpia-zambrow-victoria-sulejowek/dvQnuCGu/'

# Row 102
$ws.Cells.Item(102, 6).Value = 'Jagiellonia II'
$ws.Cells.Item(102, 7).Value = 1
$ws.Cells.Item(102, 8).Value = 'T. Mazowiecki'
$ws.Cells.Item(102, 9).Value = 1
$ws.Cells.Item(102, 10).Value = 2.55
$ws.Cells.Item(102, 11).Value = '19/10/2023 23:13'
$ws.Cells.Item(102, 12).Value = 2.12
$ws.Cells.Item(102, 13).Value = '21/10/2023 11:59'
$ws.Cells.Item(102, 14).Value = 3.32
$ws.Cells.Item(102, 15).Value = '19/10/2023 23:13'
$ws.Cells.Item(102, 16).Value = 3.53
$ws.Cells.Item(102, 17).Value = '21/10/2023 11:59'
$ws.Cells.Item(102, 18).Value = 2.2
$ws.Cells.Item(102, 19).Value = '19/10/2023 23:13'
$ws.Cells.Item(102, 20).Value = 2.86
$ws.Cells.Item(102, 21).Value = '21/10/2023 11:59'
$ws.Cells.Item(102, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/jagiellonia-tomaszow-mazowiecki/juflmREG/'

# Row 103
$ws.Cells.Item(103, 6).Value = 'Legia II'
$ws.Cells.Item(103, 7).Value = 2
$ws.Cells.Item(103, 8).Value = 'Warta Sieradz'
$ws.Cells.Item(103, 9).Value = 1
$ws.Cells.Item(103, 10).Value = 1.29
$ws.Cells.Item(103, 11).Value = '19/10/2023 23:13'
$ws.Cells.Item(103, 12).Value = 1.36
$ws.Cells.Item(103, 13).Value = '21/10/2023 11:50'
$ws.Cells.Item(103, 14).Value = 4.75
$ws.Cells.Item(103, 15).Value = '19/10/2023 23:13'
$ws.Cells.Item(103, 16).Value = 5
$ws.Cells.Item(103, 17).Value = '21/10/2023 11:50'
$ws.Cells.Item(103, 18).Value = 6.08
$ws.Cells.Item(103, 19).Value = '19/10/2023 23:13'
$ws.Cells.Item(103, 20).Value = 5.75
$ws.Cells.Item(103, 21).Value = '21/10/2023 11:50'
$ws.Cells.Item(103, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/legia-warta-sieradz/4nzPtPip/'

# Row 104
$ws.Cells.Item(104, 6).Value = 'Swit Mazowiecki'
$ws.Cells.Item(104, 7).Value = 5
$ws.Cells.Item(104, 8).Value = 'Pilica Bialobrzegi'
$ws.Cells.Item(104, 9).Value = 2
$ws.Cells.Item(104, 10).Value = 1.27
$ws.Cells.Item(104, 11).Value = '19/10/2023 23:13'
$ws.Cells.Item(104, 12).Value = 1.18
$ws.Cells.Item(104, 13).Value = '21/10/2023 11:58'
$ws.Cells.Item(104, 14).Value = 5
$ws.Cells.Item(104, 15).Value = '19/10/2023 23:13'
$ws.Cells.Item(104, 16).Value = 6.2
$ws.Cells.Item(104, 17).Value = '21/10/2023 11:58'
$ws.Cells.Item(104, 18).Value = 6.12
$ws.Cells.Item(104, 19).Value = '19/10/2023 23:13'
$ws.Cells.Item(104, 20).Value = 9.720000000000001
$ws.Cells.Item(104, 21).Value = '21/10/2023 11:58'
$ws.Cells.Item(104, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/swit-mazowiecki-pilica-bialobrzegi/WbFN1QT9/'

# Row 106
$ws.Cells.Item(106, 6).Value = 'Concordia Elblag'
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = 'Grodzisk M.'
$ws.Cells.Item(106, 9).Value = 4
$ws.Cells.Item(106, 10).Value = 4.13
$ws.Cells.Item(106, 11).Value = '20/10/2023 02:12'
$ws.Cells.Item(106, 12).Value = 4.49
$ws.Cells.Item(106, 13).Value = '21/10/2023 14:05'
$ws.Cells.Item(106, 14).Value = 3.86
$ws.Cells.Item(106, 15).Value = '20/10/2023 02:12'
$ws.Cells.Item(106, 16).Value = 4.07
$ws.Cells.Item(106, 17).Value = '21/10/2023 14:05'
$ws.Cells.Item(106, 18).Value = 1.54
$ws.Cells.Item(106, 19).Value = '20/10/2023 02:12'
$ws.Cells.Item(106, 20).Value = 1.57
$ws.Cells.Item(106, 21).Value = '21/10/2023 14:05'
$ws.Cells.Item(106, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/concordia-elblag-gks-pogon-grodzisk-mazowiecki/boHF3nac/'

# Row 107
$ws.Cells.Item(107, 6).Value = 'Zambrow'
$ws.Cells.Item(107, 7).Value = 4
$ws.Cells.Item(107, 8).Value = 'Wikielec'
$ws.Cells.Item(107, 9).Value = 0
$ws.Cells.Item(107, 10).Value = 2.28
$ws.Cells.Item(107, 11).Value = '20/10/2023 02:12'
$ws.Cells.Item(107, 12).Value = 2.48
$ws.Cells.Item(107, 13).Value = '21/10/2023 14:41'
$ws.Cells.Item(107, 14).Value = 3.2
$ws.Cells.Item(107, 15).Value = '20/10/2023 02:12'
$ws.Cells.Item(107, 16).Value = 3.29
$ws.Cells.Item(107, 17).Value = '21/10/2023 14:41'
$ws.Cells.Item(107, 18).Value = 2.53
$ws.Cells.Item(107, 19).Value = '20/10/2023 02:12'
$ws.Cells.Item(107, 20).Value = 2.51
$ws.Cells.Item(107, 21).Value = '21/10/2023 14:41'
$ws.Cells.Item(107, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/olimpia-zambrow-gks-wikielec/Cfrco5qT/'

# --- Append new rows 121-127 ---
$ws.Range("A120:V120").Copy()
$ws.Range("A121:V127").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 121
$ws.Cells.Item(121, 1).Value = 120
$ws.Cells.Item(121, 2).Value = 'poland'
$ws.Cells.Item(121, 3).Value = 'iii-liga-group-i'
$ws.Cells.Item(121, 4).Value = '2023-2024'
$ws.Cells.Item(121, 5).Value = 45234.5
$ws.Cells.Item(121, 6).Value = 'Legia II'
$ws.Cells.Item(121, 7).Value = 2
$ws.Cells.Item(121, 8).Value = 'Legionowo'
$ws.Cells.Item(121, 9).Value = 1
$ws.Cells.Item(121, 10).Value = 1.45
$ws.Cells.Item(121, 11).Value = '03/11/2023 00:12'
$ws.Cells.Item(121, 12).Value = 1.39
$ws.Cells.Item(121, 13).Value = '04/11/2023 11:05'
$ws.Cells.Item(121, 14).Value = 4.08
$ws.Cells.Item(121, 15).Value = '03/11/2023 00:12'
$ws.Cells.Item(121, 16).Value = 4.63
$ws.Cells.Item(121, 17).Value = '04/11/2023 11:32'
$ws.Cells.Item(121, 18).Value = 4.58
$ws.Cells.Item(121, 19).Value = '03/11/2023 00:12'
$ws.Cells.Item(121, 20).Value = 5.86
$ws.Cells.Item(121, 21).Value = '04/11/2023 11:06'
$ws.Cells.Item(121, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/legia-legionowo/xIg6CsrS/'

# Row 122
$ws.Cells.Item(122, 1).Value = 121
$ws.Cells.Item(122, 2).Value = 'poland'
$ws.Cells.Item(122, 3).Value = 'iii-liga-group-i'
$ws.Cells.Item(122, 4).Value = '2023-2024'
$ws.Cells.Item(122, 5).Value = 45234.54166666666
$ws.Cells.Item(122, 6).Value = 'Zambrow'
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 'Jagiellonia II'
$ws.Cells.Item(122, 9).Value = 1
$ws.Cells.Item(122, 10).Value = 1.89
$ws.Cells.Item(122, 11).Value = '03/11/2023 01:12'
$ws.Cells.Item(122, 12).Value = 2.36
$ws.Cells.Item(122, 13).Value = '04/11/2023 12:47'
$ws.Cells.Item(122, 14).Value = 3.5
$ws.Cells.Item(122, 15).Value = '03/11/2023 01:12'
$ws.Cells.Item(122, 16).Value = 3.63
$ws.Cells.Item(122, 17).Value = '04/11/2023 12:38'
$ws.Cells.Item(122, 18).Value = 3
$ws.Cells.Item(122, 19).Value = '03/11/2023 01:12'
$ws.Cells.Item(122, 20).Value = 2.46
$ws.Cells.Item(122, 21).Value = '04/11/2023 12:47'
$ws.Cells.Item(122, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/olimpia-zambrow-jagiellonia/bRf2DNSL/'

# Row 123
$ws.Cells.Item(123, 1).Value = 122
$ws.Cells.Item(123, 2).Value = 'poland'
$ws.Cells.Item(123, 3).Value = 'iii-liga-group-i'
$ws.Cells.Item(123, 4).Value = '2023-2024'
$ws.Cells.Item(123, 5).Value = 45234.54166666666
$ws.Cells.Item(123, 6).Value = 'Swit Mazowiecki'
$ws.Cells.Item(123, 7).Value = 1
$ws.Cells.Item(123, 8).Value = 'Grodzisk M.'
$ws.Cells.Item(123, 9).Value = 1
$ws.Cells.Item(123, 10).Value = 2.1
$ws.Cells.Item(123, 11).Value = '03/11/2023 01:12'
$ws.Cells.Item(123, 12).Value = 2.28
$ws.Cells.Item(123, 13).Value = '04/11/2023 12:58'
$ws.Cells.Item(123, 14).Value = 3.33
$ws.Cells.Item(123, 15).Value = '03/11/2023 01:12'
$ws.Cells.Item(123, 16).Value = 3.45
$ws.Cells.Item(123, 17).Value = '04/11/2023 12:58'
$ws.Cells.Item(123, 18).Value = 2.69
$ws.Cells.Item(123, 19).Value = '03/11/2023 01:12'
$ws.Cells.Item(123, 20).Value = 2.65
$ws.Cells.Item(123, 21).Value = '04/11/2023 12:58'
$ws.Cells.Item(123, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/swit-mazowiecki-gks-pogon-grodzisk-mazowiecki/6XHTTrD2/'

# Row 124
$ws.Cells.Item(124, 1).Value = 123
$ws.Cells.Item(124, 2).Value = 'poland'
$ws.Cells.Item(124, 3).Value = 'iii-liga-group-i'
$ws.Cells.Item(124, 4).Value = '2023-2024'
$ws.Cells.Item(124, 5).Value = 45234.58333333334
$ws.Cells.Item(124, 6).Value = 'Mlawa'
$ws.Cells.Item(124, 7).Value = 2
$ws.Cells.Item(124, 8).Value = 'LKS Lomza'
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 1.72
$ws.Cells.Item(124, 11).Value = '03/11/2023 02:13'
$ws.Cells.Item(124, 12).Value = 1.83
$ws.Cells.Item(124, 13).Value = '04/11/2023 13:52'
$ws.Cells.Item(124, 14).Value = 3.74
$ws.Cells.Item(124, 15).Value = '03/11/2023 02:13'
$ws.Cells.Item(124, 16).Value = 4.01
$ws.Cells.Item(124, 17).Value = '04/11/2023 13:52'
$ws.Cells.Item(124, 18).Value = 3.33
$ws.Cells.Item(124, 19).Value = '03/11/2023 02:13'
$ws.Cells.Item(124, 20).Value = 3.25
$ws.Cells.Item(124, 21).Value = '04/11/2023 13:05'
$ws.Cells.Item(124, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/mks-mlawa-lks-lomza/4Swo5az2/'

# Row 125
$ws.Cells.Item(125, 1).Value = 124
$ws.Cells.Item(125, 2).Value = 'poland'
$ws.Cells.Item(125, 3).Value = 'iii-liga-group-i'
$ws.Cells.Item(125, 4).Value = '2023-2024'
$ws.Cells.Item(125, 5).Value = 45234.58333333334
$ws.Cells.Item(125, 6).Value = 'Concordia Elblag'
$ws.Cells.Item(125, 7).Value = 4
$ws.Cells.Item(125, 8).Value = 'Sulejowek'
$ws.Cells.Item(125, 9).Value = 1
$ws.Cells.Item(125, 10).Value = 2.67
$ws.Cells.Item(125, 11).Value = '03/11/2023 02:13'
$ws.Cells.Item(125, 12).Value = 2.81
$ws.Cells.Item(125, 13).Value = '04/11/2023 13:52'
$ws.Cells.Item(125, 14).Value = 3.23
$ws.Cells.Item(125, 15).Value = '03/11/2023 02:13'
$ws.Cells.Item(125, 16).Value = 3.35
$ws.Cells.Item(125, 17).Value = '04/11/2023 13:52'
$ws.Cells.Item(125, 18).Value = 2.16
$ws.Cells.Item(125, 19).Value = '03/11/2023 02:13'
$ws.Cells.Item(125, 20).Value = 2.22
$ws.Cells.Item(125, 21).Value = '04/11/2023 13:52'
$ws.Cells.Item(125, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/concordia-elblag-victoria-sulejowek/jPZt6uLe/'

# Row 126
$ws.Cells.Item(126, 1).Value = 125
$ws.Cells.Item(126, 2).Value = 'poland'
$ws.Cells.Item(126, 3).Value = 'iii-liga-group-i'
$ws.Cells.Item(126, 4).Value = '2023-2024'
$ws.Cells.Item(126, 5).Value = 45234.58333333334
$ws.Cells.Item(126, 6).Value = 'GKS Belchatow'
$ws.Cells.Item(126, 7).Value = 2
$ws.Cells.Item(126, 8).Value = 'T. Mazowiecki'
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 1.93
$ws.Cells.Item(126, 11).Value = '03/11/2023 02:13'
$ws.Cells.Item(126, 12).Value = 1.77
$ws.Cells.Item(126, 13).Value = '04/11/2023 13:41'
$ws.Cells.Item(126, 14).Value = 3.43
$ws.Cells.Item(126, 15).Value = '03/11/2023 02:13'
$ws.Cells.Item(126, 16).Value = 3.75
$ws.Cells.Item(126, 17).Value = '04/11/2023 13:41'
$ws.Cells.Item(126, 18).Value = 2.96
$ws.Cells.Item(126, 19).Value = '03/11/2023 02:13'
$ws.Cells.Item(126, 20).Value = 3.66
$ws.Cells.Item(126, 21).Value = '04/11/2023 13:41'
$ws.Cells.Item(126, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/gks-belchatow-tomaszow-mazowiecki/44ebE3DF/'

# Row 127
$ws.Cells.Item(127, 1).Value = 126
$ws.Cells.Item(127, 2).Value = 'poland'
$ws.Cells.Item(127, 3).Value = 'iii-liga-group-i'
$ws.Cells.Item(127, 4).Value = '2023-2024'
$ws.Cells.Item(127, 5).Value = 45234.58333333334
$ws.Cells.Item(127, 6).Value = 'Pilica Bialobrzegi'
$ws.Cells.Item(127, 7).Value = 0
$ws.Cells.Item(127, 8).Value = 'Skierniewice'
$ws.Cells.Item(127, 9).Value = 4
$ws.Cells.Item(127, 10).Value = 2.88
$ws.Cells.Item(127, 11).Value = '03/11/2023 02:13'
$ws.Cells.Item(127, 12).Value = 3.2
$ws.Cells.Item(127, 13).Value = '03/11/2023 11:31'
$ws.Cells.Item(127, 14).Value = 3.39
$ws.Cells.Item(127, 15).Value = '03/11/2023 02:13'
$ws.Cells.Item(127, 16).Value = 3.56
$ws.Cells.Item(127, 17).Value = '04/11/2023 12:02'
$ws.Cells.Item(127, 18).Value = 1.97
$ws.Cells.Item(127, 19).Value = '03/11/2023 02:13'
$ws.Cells.Item(127, 20).Value = 1.93
$ws.Cells.Item(127, 21).Value = '03/11/2023 11:31'
$ws.Cells.Item(127, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-i/pilica-bialobrzegi-unia-skierniewice/pOGXS2S8/'

